# Insert a new data row just before the current row 787 (Excel's Rows.Insert
# shifts row 787 and everything below it down by one, growing the used range
# from A1:R842 to A1:R843), then populate the newly-inserted row 787 with the
# new price-sheet entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(787).Insert()

$ws.Cells.Item(787, 1).Value = 8
$ws.Cells.Item(787, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(787, 3).Value = "Coquimbo"
$ws.Cells.Item(787, 4).Value = 44931
$ws.Cells.Item(787, 5).Value = 4
$ws.Cells.Item(787, 6).Value = 100112045
$ws.Cells.Item(787, 7).Value = "Zapallo"
$ws.Cells.Item(787, 8).Value = "Camote"
$ws.Cells.Item(787, 9).Value = "1a (cosecha)"
$ws.Cells.Item(787, 10).Value = 1660
$ws.Cells.Item(787, 11).Value = 900
$ws.Cells.Item(787, 12).Value = 1000
$ws.Cells.Item(787, 13).Value = 950
$ws.Cells.Item(787, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(787, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(787, 16).Value = 950
$ws.Cells.Item(787, 17).Value = 1
$ws.Cells.Item(787, 18).Value = "Hortaliza"
